# Scheduled-runner market data refresh: update computed price/profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to the
# latest pulled values. Plain value overwrites (no formulas in this
# workbook) - a handful of rows also lose/gain a cell in M or N where the
# source feed stopped/started reporting that figure, so those are cleared
# or set explicitly to match.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2070.3572
$ws.Range("I18").Value = 1207.25
$ws.Range("K18").Value = 1207.25
$ws.Range("M18").Value = -923.25

$ws.Range("H33").Value = 524.0909
$ws.Range("J33").Value = 329.25
$ws.Range("L33").Value = 329.25
$ws.Range("N33").Value = -787.25

$ws.Range("H43").Value = 11842.143
$ws.Range("J43").Value = 9649.5
$ws.Range("L43").Value = 9649.5
$ws.Range("N43").Value = -9787.5

$ws.Range("H48").Value = 6199.2
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H56").Value = 6199.2
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H132").Value = 4168.0684
$ws.Range("I132").Value = 4377.0312
$ws.Range("J132").Value = 3610.8333
$ws.Range("K132").Value = 13131.0936
$ws.Range("L132").Value = 10832.4999
$ws.Range("M132").Value = -10601.0936
$ws.Range("N132").Value = -15892.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1972.8734
$ws.Range("I32").Value = 1858.9305
$ws.Range("K32").Value = 1858.9305
$ws.Range("M32").Value = -1571.9305

$ws.Range("H122").Value = 4951.4
$ws.Range("I122").Value = 4951.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14854.2
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 6067.979
$ws.Range("I132").Value = 1879.579
$ws.Range("K132").Value = 5638.737
$ws.Range("M132").Value = -3108.737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 73577.8
$ws.Range("I135").Value = 61000
$ws.Range("K135").Value = 61000
$ws.Range("M135").Value = -55930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1550
$ws.Range("I19").Value = 150
$ws.Range("K19").Value = 150
$ws.Range("M19").Value = 20

$ws.Range("H24").Value = 1550
$ws.Range("I24").Value = 150
$ws.Range("K24").Value = 150
$ws.Range("M24").Value = 20

$ws.Range("H94").Value = 2023
$ws.Range("I94").Value = 2203.6667
$ws.Range("J94").Value = 1962.7778
$ws.Range("K94").Value = 2203.6667
$ws.Range("L94").Value = 1962.7778
$ws.Range("M94").Value = -1752.6667
$ws.Range("N94").Value = -2864.7778

$ws.Range("H122").Value = 5081.1
$ws.Range("I122").Value = 4969.143
$ws.Range("K122").Value = 14907.429
$ws.Range("M122").Value = -12457.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 76923300
$ws.Range("I23").Value = 140.75
$ws.Range("J23").Value = 111111370
$ws.Range("K23").Value = 422.25
$ws.Range("L23").Value = 333334110
$ws.Range("M23").Value = -187.25
$ws.Range("N23").Value = -333334580

$ws.Range("H39").Value = 2670.5
$ws.Range("J39").Value = 2726.818
$ws.Range("L39").Value = 8180.454000000001
$ws.Range("N39").Value = -8768.454000000002

$ws.Range("H55").Value = 10402.765
$ws.Range("J55").Value = 11643.333
$ws.Range("L55").Value = 34929.999
$ws.Range("N55").Value = -35283.999

$ws.Range("H75").Value = 3531.3
$ws.Range("J75").Value = 5506
$ws.Range("L75").Value = 16518
$ws.Range("N75").Value = -18514

$ws.Range("H78").Value = 3531.3
$ws.Range("J78").Value = 5506
$ws.Range("L78").Value = 49554
$ws.Range("N78").Value = -59538

$ws.Range("H113").Value = 1225.3125
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1225.3125
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8015.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 14795.054
$ws.Range("I102").Value = 1335.129
$ws.Range("J102").Value = 84338
$ws.Range("K102").Value = 1335.129
$ws.Range("L102").Value = 84338
$ws.Range("M102").Value = 286.8710000000001
$ws.Range("N102").Value = -87582

$ws.Range("H111").Value = 51159.668
$ws.Range("J111").Value = 51159.668
$ws.Range("L111").Value = 51159.668
$ws.Range("N111").Value = -57293.668

$ws.Range("H113").Value = 4699.3335
$ws.Range("I113").Value = 8997
$ws.Range("J113").Value = 2550.5
$ws.Range("K113").Value = 8997
$ws.Range("L113").Value = 2550.5
$ws.Range("M113").Value = -6827
$ws.Range("N113").Value = -6890.5

$ws.Range("H132").Value = 3649.9832
$ws.Range("J132").Value = 2338.4666
$ws.Range("L132").Value = 7015.399800000001
$ws.Range("N132").Value = -12075.3998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5871.25
$ws.Range("I7").Value = 5496
$ws.Range("K7").Value = 5496
$ws.Range("M7").Value = -5384

$ws.Range("H22").Value = 2977.4285
$ws.Range("I22").Value = 911.6667
$ws.Range("J22").Value = 3540.818
$ws.Range("K22").Value = 911.6667
$ws.Range("L22").Value = 3540.818
$ws.Range("M22").Value = -616.6667
$ws.Range("N22").Value = -4130.818

$ws.Range("H27").Value = 2977.4285
$ws.Range("I27").Value = 911.6667
$ws.Range("J27").Value = 3540.818
$ws.Range("K27").Value = 911.6667
$ws.Range("L27").Value = 3540.818
$ws.Range("M27").Value = -804.6667
$ws.Range("N27").Value = -3754.818

$ws.Range("H40").Value = 9338.132
$ws.Range("I40").Value = 9081.044
$ws.Range("J40").Value = 9732.333000000001
$ws.Range("K40").Value = 9081.044
$ws.Range("L40").Value = 9732.333000000001
$ws.Range("M40").Value = -8945.044
$ws.Range("N40").Value = -10004.333

$ws.Range("H46").Value = 12528.777
$ws.Range("J46").Value = 15857.714
$ws.Range("L46").Value = 15857.714
$ws.Range("N46").Value = -16233.714

$ws.Range("H110").Value = 45000
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 7327.408
$ws.Range("I122").Value = 6503.476
$ws.Range("J122").Value = 12271
$ws.Range("K122").Value = 19510.428
$ws.Range("L122").Value = 36813
$ws.Range("M122").Value = -17060.428
$ws.Range("N122").Value = -41713

$ws.Range("H126").Value = 5871.25
$ws.Range("I126").Value = 5496
$ws.Range("K126").Value = 16488
$ws.Range("M126").Value = -14018

$ws.Range("H132").Value = 5392.3335
$ws.Range("I132").Value = 3884.75
$ws.Range("K132").Value = 11654.25
$ws.Range("M132").Value = -9124.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3409.1428
$ws.Range("I96").Value = 3193.111
$ws.Range("J96").Value = 3798
$ws.Range("K96").Value = 3193.111
$ws.Range("L96").Value = 3798
$ws.Range("M96").Value = -1820.111
$ws.Range("N96").Value = -6544

$ws.Range("H122").Value = 4624.222
$ws.Range("I122").Value = 4736.067
$ws.Range("K122").Value = 14208.201
$ws.Range("M122").Value = -11758.201

$ws.Range("H126").Value = 7592.1035
$ws.Range("I126").Value = 8967.143
$ws.Range("J126").Value = 6308.7334
$ws.Range("K126").Value = 26901.429
$ws.Range("L126").Value = 18926.2002
$ws.Range("M126").Value = -24431.429
$ws.Range("N126").Value = -23866.2002

$ws.Range("H132").Value = 1800.5
$ws.Range("I132").Value = 1579.8
$ws.Range("K132").Value = 4739.4
$ws.Range("M132").Value = -2209.4

$ws.Range("H133").Value = 121857
$ws.Range("J133").Value = 121857
$ws.Range("L133").Value = 121857
$ws.Range("N133").Value = -131977

$ws.Range("H136").Value = 6104.92
$ws.Range("I136").Value = 6286.2163
$ws.Range("K136").Value = 18858.6489
$ws.Range("M136").Value = -16308.6489
